$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.822.35'
$ws.Range("E2").Value = '  -1.43%  '

$ws.Range("D3").Value = '3.067.34'
$ws.Range("E3").Value = '  -0.64%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '535.90'
$ws.Range("E5").Value = '  -3.30%  '

$ws.Range("D6").Value = '133.39'
$ws.Range("E6").Value = '  -2.83%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = '3.060.36'
$ws.Range("E8").Value = '  -0.50%  '

$ws.Range("E9").Value = '  +0.16%  '

$ws.Range("E10").Value = '  -1.15%  '

$ws.Range("D11").Value = '6.15'
$ws.Range("E11").Value = '  -7.46%  '

$ws.Range("D12").Value = '0.452'
$ws.Range("E12").Value = '  +0.67%  '

$ws.Range("E13").Value = '  +3.55%  '

$ws.Range("E14").Value = '  -2.37%  '

$ws.Range("D15").Value = '3.557.54'
$ws.Range("E15").Value = '  -0.65%  '

$ws.Range("D16").Value = '62.879.85'
$ws.Range("E16").Value = '  -1.31%  '

$ws.Range("E17").Value = '  -0.36%  '

$ws.Range("D18").Value = '3.062.35'
$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("E19").Value = '  +0.53%  '

$ws.Range("D20").Value = '480.20'

$ws.Range("D21").Value = '13.29'
$ws.Range("E21").Value = '  -1.58%  '

$ws.Range("D22").Value = '0.692'
$ws.Range("E22").Value = '  -1.00%  '

$ws.Range("E23").Value = '  -0.85%  '

$ws.Range("D24").Value = '78.86'
$ws.Range("E24").Value = '  +2.40%  '

$ws.Range("E25").Value = '  -1.73%  '

$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("E27").Value = '  -2.09%  '

$ws.Range("E28").Value = '  -2.05%  '

$ws.Range("E29").Value = '  -0.16%  '

$ws.Range("E30").Value = '  -0.30%  '

$ws.Range("E31").Value = '  -8.29%  '

$ws.Range("E32").Value = '  +0.80%  '

$ws.Range("E33").Value = '  -6.91%  '

$ws.Range("D34").Value = '56.62'
$ws.Range("E34").Value = '  -2.12%  '

$ws.Range("E35").Value = '  +4.12%  '

$ws.Range("E36").Value = '  +1.99%  '

$ws.Range("D37").Value = '476.40'
$ws.Range("E37").Value = '  -10.18%  '

$ws.Range("D38").Value = '0.0394'
$ws.Range("E38").Value = '  -3.93%  '

$ws.Range("D39").Value = '3.083.88'
$ws.Range("E39").Value = '  +1.24%  '

$ws.Range("E40").Value = '  +0.49%  '

$ws.Range("E41").Value = '  -1.45%  '

$ws.Range("D42").Value = '8.09'
$ws.Range("E42").Value = '  +0.61%  '

$ws.Range("D43").Value = '2.64'
$ws.Range("E43").Value = '  +3.32%  '

$ws.Range("E44").Value = '  +0.92%  '

$ws.Range("E46").Value = '  +9.80%  '

$ws.Range("D47").Value = '121.19'
$ws.Range("E47").Value = '  -0.74%  '

$ws.Range("D48").Value = '2.01'
$ws.Range("E48").Value = '  -2.08%  '

$ws.Range("E49").Value = '  +2.17%  '

$ws.Range("E50").Value = '  +1.76%  '

$ws.Range("E51").Value = '  +1.62%  '
